# Updating the models for NRG, PCSun and Ulmeni
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..33 (column A = value, column B = timestamp serial)
$data = @(
    @(5801, 46068.95833333334),
    @(5817, 46068.96875),
    @(5856, 46068.97916666666),
    @(5625, 46068.98958333334),
    @(5474, 46069),
    @(5406, 46069.01041666666),
    @(5553, 46069.02083333334),
    @(5551, 46069.03125),
    @(5548, 46069.04166666666),
    @(5535, 46069.05208333334),
    @(5550, 46069.0625),
    @(5508, 46069.08333333334),
    @(5498, 46069.09375),
    @(5484, 46069.10416666666),
    @(5523, 46069.11458333334),
    @(5654, 46069.125),
    @(5662, 46069.13541666666),
    @(5738, 46069.14583333334),
    @(5803, 46069.15625),
    @(5918, 46069.16666666666),
    @(6001, 46069.17708333334),
    @(6046, 46069.1875),
    @(6180, 46069.19791666666),
    @(6202, 46069.20833333334),
    @(6377, 46069.21875),
    @(6566, 46069.22916666666),
    @(6777, 46069.23958333334),
    @(7079, 46069.25),
    @(7220, 46069.26041666666),
    @(7403, 46069.27083333334),
    @(7623, 46069.28125),
    @(7848, 46069.29166666666)
)

# Old sheet had data through row 53; new sheet only goes to row 33.
# Clear out the old trailing rows (34..53) completely.
$ws.Range("A34:B53").Clear()

# Write the new values into rows 2..33
$rowIndex = 2
foreach ($pair in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $pair[0]
    $ws.Cells.Item($rowIndex, 2).Value = $pair[1]
    $rowIndex++
}
